$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.824.33'
$ws.Range("D3").Value = '3.249.12'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''583.81'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").Value = '''182.88'
$ws.Range("E6").Value = '  +3.94%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''0.598'
$ws.Range("E8").Value = '  -1.17%  '
$ws.Range("E9").Value = '  +4.71%  '
$ws.Range("D10").Value = '''6.68'
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("E11").Value = '  +2.46%  '
$ws.Range("D12").Value = '3.811.53'
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D14").Value = '''28.65'
$ws.Range("E14").Value = '  +3.30%  '
$ws.Range("D15").Value = '67.821.20'
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("D16").Value = '''0.0000171'
$ws.Range("E16").Value = '  +2.44%  '
$ws.Range("D17").Value = '3.254.04'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("E18").Value = '  +0.75%  '
$ws.Range("E19").Value = '  +2.43%  '
$ws.Range("D20").Value = '''379.89'
$ws.Range("E20").Value = '  +3.49%  '
$ws.Range("D21").Value = '''7.63'
$ws.Range("E21").Value = '  +2.66%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '''71.26'
$ws.Range("E23").Value = '  +1.76%  '
$ws.Range("D24").Value = '''0.512'
$ws.Range("E24").Value = '  +1.63%  '
$ws.Range("E25").Value = '  +1.20%  '
$ws.Range("D26").Value = '''9.91'
$ws.Range("E26").Value = '  +2.22%  '
$ws.Range("E27").Value = '  +2.10%  '
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D30").Value = '''5.66'
$ws.Range("E30").Value = '  +1.64%  '
$ws.Range("D31").Value = '''22.84'
$ws.Range("E31").Value = '  +2.09%  '
$ws.Range("D32").Value = '''7.08'
$ws.Range("E32").Value = '  +5.05%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +3.09%  '
$ws.Range("E35").Value = '  +4.54%  '
$ws.Range("D36").Value = '''162.07'
$ws.Range("E36").Value = '  -6.72%  '
$ws.Range("D37").Value = '''0.837'
$ws.Range("E37").Value = '  -1.23%  '
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''6.70'
$ws.Range("E39").Value = '  +5.34%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '''26.47'
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").Value = '''4.57'
$ws.Range("E41").Value = '  +7.07%  '
$ws.Range("D42").Value = '''2.58'
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '''41.20'
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''25.43'
$ws.Range("E44").Value = '  +4.01%  '
$ws.Range("D45").Value = '''345.55'
$ws.Range("E45").Value = '  +3.54%  '
$ws.Range("D46").Value = '''0.0685'
$ws.Range("E46").Value = '  +2.61%  '
$ws.Range("D47").Value = '2.621.97'
$ws.Range("E47").Value = '  -2.99%  '
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("D50").Value = '''0.993'
$ws.Range("E50").Value = '  +1.80%  '
$ws.Range("E51").Value = '  +2.85%  '

# Reset style on cells written with a leading quote-prefix (forced text)
# so they do not pick up an explicit "quote prefix" cell style different
# from the workbook default (matches original formatting: no style override).
$resetCells = @("D5","D6","D8","D10","D14","D16","D20","D21","D23","D24","D26","D30","D31","D32","D36","D37","D39","D40","D41","D42","D43","D44","D45","D46","D50")
foreach ($addr in $resetCells) {
    $ws.Range($addr).Style = "Normal"
}
